$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.011.16'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.299.85'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.13'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.10'
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.21'
$ws.Range('E10').Value = '  +3.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('E12').Value = '  +2.73%  '
$ws.Range('E13').Value = '  +4.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.98'
$ws.Range('E14').Value = '  +16.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.77'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.657.26'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.288.60'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('E18').Value = '  +4.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.921.90'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.35'
$ws.Range('E20').Value = '  +8.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.10'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.88'
$ws.Range('E23').Value = '  +1.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.52'
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('E25').Value = '  +13.54%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.46'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.70'
$ws.Range('E28').Value = '  +4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '168.81'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.69'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.15'
$ws.Range('E32').Value = '  +1.18%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.55'
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  +3.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.89'
$ws.Range('E37').Value = '  +4.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0693'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.102'
$ws.Range('E39').Value = '  +3.23%  '
$ws.Range('E40').Value = '  +4.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.79'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.994.60'
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('E46').Value = '  +4.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.55'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.85'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '56.48'
$ws.Range('E49').Value = '  +8.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.528.71'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('E51').Value = '  +3.23%  '
